$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsSchedule.Range("E2").Value = 382.7884124999999
$wsSchedule.Range("F2").Value = 8.43889798280423
$wsSchedule.Range("E3").Value = 446.5496879999999
$wsSchedule.Range("F3").Value = 29.53370952380952
$wsDetailed.Range("B5").Value = 66.29066
$wsDetailed.Range("B6").Value = 66.0025
$wsDetailed.Range("B7").Value = 66.91710999999999
$wsDetailed.Range("C7").Value = "historical"
$wsDetailed.Range("C8").Value = "historical"
$wsDetailed.Range("C9").Value = "historical"
$wsDetailed.Range("B10").Value = 77.94
$wsDetailed.Range("B11").Value = 78
$wsDetailed.Range("B12").Value = 77.94
$wsDetailed.Range("B13").Value = 80.5187
$wsDetailed.Range("B14").Value = 76.62188999999999
$wsDetailed.Range("B15").Value = 59.02779
$wsDetailed.Range("B16").Value = 35.88
$wsDetailed.Range("B17").Value = 0.51
$wsDetailed.Range("B18").Value = -5.50985
$wsDetailed.Range("B19").Value = -6.22853
$wsDetailed.Range("B20").Value = -6.49292
$wsDetailed.Range("B21").Value = -7.8034
$wsDetailed.Range("B22").Value = -7.14347
$wsDetailed.Range("B23").Value = -7.53039
$wsDetailed.Range("B24").Value = -7.4391
$wsDetailed.Range("B25").Value = -5.58973
$wsDetailed.Range("B26").Value = -3.64783
$wsDetailed.Range("B27").Value = -6.49292
$wsDetailed.Range("B28").Value = -6.17453
$wsDetailed.Range("B29").Value = -5.89628
$wsDetailed.Range("B31").Value = -0.87608
$wsDetailed.Range("B33").Value = -12.01
$wsDetailed.Range("B34").Value = -10
$wsDetailed.Range("B35").Value = -10
$wsDetailed.Range("B36").Value = -11.01
$wsDetailed.Range("B37").Value = -7.32903
$wsDetailed.Range("B38").Value = -1.1176
$wsDetailed.Range("B39").Value = 7.27165
$wsDetailed.Range("B40").Value = 29.54996
$wsDetailed.Range("B42").Value = 55.33037
$wsDetailed.Range("B43").Value = 61.94424
$wsDetailed.Range("B44").Value = 57.04922
$wsDetailed.Range("B45").Value = 57.03883
$wsDetailed.Range("B46").Value = 55.69702
$wsDetailed.Range("B49").Value = 56.98
